$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from the "Uchwala Czlonkow..." paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Locate the last paragraph in the body -- "Przewodniczacy zebrania ...
#    Sekretarz zebrania" -- and replace it wholesale (content + mark) with
#    the new block of paragraphs described by the diff: the reformatted
#    Przewodniczacy/Sekretarz line, a blank line, and the names line.
$targetPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$targetRange = $targetPara.Range

$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="-5" w:hanging="10"/><w:rPr><w:rFonts w:ascii="Liberation Serif" w:eastAsia="DejaVu Sans" w:hAnsi="Liberation Serif" w:cs="DejaVu Sans"/><w:color w:val="auto"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Przewodniczący zebrania </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Sekretarz zebrania</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="-5" w:hanging="10"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="-5" w:hanging="10"/></w:pPr><w:r><w:t xml:space="preserve">       Wojciech Zając                                                                                                              Rafał Korzeniewski</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($xmlFragment)

# 3. InsertXML leaves one more paragraph behind, carrying the ORIGINAL
#    paragraph's mark/formatting (spacing + ind) but with its runs cleared
#    out -- exactly the shell the final (bookmark-only) paragraph in the
#    diff needs. Put the bookmark back into that paragraph. (A fresh
#    Range(start,start) is used -- not a collapsed Paragraph.Range -- so
#    the bookmark start/end land together instead of drifting.)
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$pos = $finalPara.Range.Start
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
